# Update gh-pages output data for 南宁-漫展信息 workbook.
#
# Sheet "展览" (exhibitions):
#   row 2: 想去人数(F) 7264 -> 7270 ; 最低票价(G) "已售罄" -> 62
#   row 4: 想去人数(F) 127 -> 129
#   row 5: 想去人数(F) 179 -> 181
#   row 7: 想去人数(F) 97 -> 100
#   row 8: 想去人数(F) 615 -> 619
#
# Sheet "演出" (performances):
#   row 2: 最低票价(G) 100 -> "不可售"
#
# Sheet "全部类型" (all types, combines the two sheets above) mirrors the
# same changes on its own rows.

$wb = $excel.ActiveWorkbook

$shows = $wb.Worksheets.Item("展览")
$shows.Range("F2").Value = 7270
$shows.Range("G2").Value = 62
$shows.Range("F4").Value = 129
$shows.Range("F5").Value = 181
$shows.Range("F7").Value = 100
$shows.Range("F8").Value = 619

$perf = $wb.Worksheets.Item("演出")
$perf.Range("G2").Value = "不可售"

$all = $wb.Worksheets.Item("全部类型")
$all.Range("F2").Value = 7270
$all.Range("G2").Value = 62
$all.Range("G4").Value = "不可售"
$all.Range("F5").Value = 129
$all.Range("F6").Value = 181
$all.Range("F9").Value = 100
$all.Range("F10").Value = 619
